# Update the cached "datetimeFigureOut" field text on the slide master
# and every slide layout from 2019/2/15 to 2019/2/16.

$p = $ppt.ActivePresentation

function Update-DateField {
    param($shapes)
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq "2019/2/15") {
                $tr.Text = "2019/2/16"
            }
        }
    }
}

# Slide master
Update-DateField $p.SlideMaster.Shapes

# Every slide layout belonging to the master
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $cl = $layouts.Item($li)
    Update-DateField $cl.Shapes
}
